$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting (border/bold/alignment) used
# by the other header cells (e.g. G1) before setting its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the new "Save" column values for the data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
